$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (ns) values
$ws.Range("C2").Value = 86
$ws.Range("C3").Value = 87
$ws.Range("C4").Value = 90
$ws.Range("C5").Value = 88
$ws.Range("C7").Value = 87
$ws.Range("C8").Value = 89
$ws.Range("C10").Value = 93
$ws.Range("C11").Value = 92
$ws.Range("C12").Value = 86
$ws.Range("C13").Value = 87
$ws.Range("C14").Value = 88
$ws.Range("C15").Value = 84
$ws.Range("C16").Value = 88
$ws.Range("C17").Value = 86
$ws.Range("C19").Value = 89
$ws.Range("C20").Value = 88
$ws.Range("C21").Value = 92
$ws.Range("C22").Value = 88
$ws.Range("C23").Value = 90
$ws.Range("C24").Value = 85
$ws.Range("C25").Value = 93
$ws.Range("C26").Value = 87
$ws.Range("C27").Value = 88
$ws.Range("C28").Value = 88
$ws.Range("C29").Value = 86
$ws.Range("C30").Value = 89
$ws.Range("C31").Value = 87
$ws.Range("C32").Value = 93
$ws.Range("C33").Value = 87
$ws.Range("C34").Value = 92
$ws.Range("C35").Value = 92
$ws.Range("C36").Value = 87
$ws.Range("C37").Value = 89
$ws.Range("C38").Value = 88
$ws.Range("C39").Value = 86
$ws.Range("C40").Value = 90
$ws.Range("C41").Value = 89
$ws.Range("C42").Value = 88
$ws.Range("C44").Value = 89
$ws.Range("C45").Value = 90

# Update column D (nr) values that were previously empty
$ws.Range("D14").Value = 89
$ws.Range("D18").Value = 90

# Update the selection in the sheet view
$ws.Range("D19").Select()
